$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4455971197233026
$ws.Range("D2").Value = 0.6602424197275307

$ws.Range("C3").Value = 0.3606132201812469
$ws.Range("D3").Value = 0.7218280105390669

$ws.Range("C4").Value = 1.31561766665408
$ws.Range("D4").Value = 0.201848726231606

$ws.Range("C5").Value = 0.2805198501366006
$ws.Range("D5").Value = 0.7816997324742947

$ws.Range("C6").Value = 0.02881861154526595
$ws.Range("D6").Value = 0.9772690910377133

$ws.Range("C7").Value = 1.156799736153584
$ws.Range("D7").Value = 0.2597611767791046

$ws.Range("C8").Value = 0.01738025307447338
$ws.Range("D8").Value = 0.986289925016917

$ws.Range("C9").Value = 0.6837731997478845
$ws.Range("D9").Value = 0.5012576399160147

$ws.Range("C10").Value = -0.01265155390170644
$ws.Range("D10").Value = 0.9900198205515329

$ws.Range("C11").Value = -0.7406901920528188
$ws.Range("D11").Value = 0.4667138577045544
